$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "[Python - 프로그래머스] 로또의 최고 순위와 최저 순위"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Python-%ED%94%84%EB%A1%9C%EA%B7%B8%EB%9E%98%EB%A8%B8%EC%8A%A4-%EB%A1%9C%EB%98%90%EC%9D%98-%EC%B5%9C%EA%B3%A0-%EC%88%9C%EC%9C%84%EC%99%80-%EC%B5%9C%EC%A0%80-%EC%88%9C%EC%9C%84"

$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

$ws.Range("D50").Value = "가장 쉬운 autoencoder(자동암호기)를 이용한 차원 축소"
$ws.Range("E50").Value = "http://incredible.egloos.com/7537539"

$ws.Range("D51").Value = "[python] 파이썬스럽게 코딩하는 법, PEP8 정리"
$ws.Range("E51").Value = "https://bskyvision.com/1249"
